$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 2 data (pout-pourri fixes) ---
$ws.Range("C2").Value = "Deep Purple In Rock (1970)"
$ws.Range("G2").Value = "Harvest - SHVL 777, Harvest - 1E 062 º 91442"
$ws.Range("I2").Value = "Speed King / Bloodsucker / Child In Time / Flight Of The Rat / Into The Fire / Living Wreck / Hard Lovin' Man"
$ws.Range("J2").Value = "Escolha o estilo"

# --- New row 3: "evento atualizar estado" form submission (teste row) ---
$ws.Range("B3").Value = "teste"
$ws.Range("C3").Value = "aksdjfbksdfbs"
$ws.Range("F3").Value = "açsndalçjkdbna"
$ws.Range("G3").Value = "sçdjnsdf"
# Force H3 to stay a text cell (matches source "Ano" column being text, not numeric)
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "1970"
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Value = "Speed King / Bloodsucker / Child In Time / Flight Of The Rat / Into The Fire / Living Wreck / Hard Lovin' Man"
$ws.Range("J3").Value = "Escolha o estilo"

# --- Rows 4-13: placeholder "None" description entries ---
for ($r = 4; $r -le 13; $r++) {
    $ws.Cells.Item($r, 9).Value = "None"
}

# --- Row 14: another partial submission ---
$ws.Range("B14").Value = "a"
$ws.Cells.Item(14, 9).Value = "None"

# --- Row 15: another placeholder row ---
$ws.Cells.Item(15, 9).Value = "None"
